# This edit re-shuffles the (Fecha, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Precio $/Kg) values across the data rows
# (rows 2-35) of the sheet. Row 16 keeps its original values. Columns
# D, J, K, L, M, P are affected; all other columns are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the values being redistributed across rows.
$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot every row's current (before) values for the affected columns,
# so the permutation can be applied safely without clobbering source data.
$before = @{}
for ($r = 2; $r -le 35; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

# Mapping of destination row -> source row whose original values it
# should receive (row 16 is unchanged, i.e. maps to itself).
$mapping = @{
    2  = 12
    3  = 11
    4  = 2
    5  = 8
    6  = 22
    7  = 31
    8  = 29
    9  = 21
    10 = 6
    11 = 25
    12 = 15
    13 = 20
    14 = 3
    15 = 18
    16 = 16
    17 = 24
    18 = 34
    19 = 14
    20 = 17
    21 = 26
    22 = 32
    23 = 33
    24 = 9
    25 = 27
    26 = 28
    27 = 30
    28 = 10
    29 = 23
    30 = 7
    31 = 19
    32 = 35
    33 = 4
    34 = 5
    35 = 13
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $before[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
